# Add an "Address" column at F, shifting the existing "District" column (F) to G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column F (District).
$ws.Columns.Item(6).Insert()

# New column header
$ws.Cells.Item(2, 6).Value = "Address"

# Row 3 sub-header stays blank for the new Address column (matches original layout)
$ws.Cells.Item(3, 6).Value = ""

# Address values extracted from the combined Name/Address text in column B
$addresses = @{
    4 = "G J C H S B H Road"
    5 = "G J C HuliyarudurgaKunigal"
    6 = "G P U C High School HonnavalliTiptur"
    7 = "G J V Thirumani Pavagada"
    8 = "G H S BukkapatnaKortagere"
    9 = "G H S DuddanahalliKoratagere"
    10 = "G J C Pavagada"
    11 = "Assistant Master G H S BaraguruC N Halli"
    12 = "Bapuji High School Bevnahalli"
    13 = "S S A H S Kunigal"
    14 = "G H S GangasagaraPavagada"
    15 = "G H S IdagurGubbi"
    16 = "G G J C (HS) Sira"
    17 = "G H S GodrahalliKoratagere"
    18 = "G H S Kyathasandra"
    19 = "Sri Kalidasa High School BidareGubbi"
    20 = "G J C Nagasandra(K G Temple) Gubbi"
    21 = "G H S AmmanaghattaGubbi"
    23 = "G J C ThirumaniPavagada"
    24 = "Sri K N Rural High School Ranganahalli Gate Gubbi"
    25 = "G G J C Tiptur"
    26 = "Sri Uma Pragathi High School Kyathsandra"
    27 = "G J C GuligenahalliSira"
    28 = "G H S VenkatapuraPavagada"
    29 = "S L H S C NandihalliGubbi"
    30 = "G J C (HS) Sira"
    31 = "Assistant Master G U H S HuliyarC N Halli"
    32 = "G J C (H S) I K ColonyKoratagere"
    33 = "S G C P U C Hebbur"
    34 = "G G J C Gubbi"
    35 = "Kalidasa Junior College"
    36 = "Sri Sai Gurukula Residential School Honnali"
    37 = "G H S YalekadakaluKunigal"
    38 = "G H S GowdanakatteTiptur"
    39 = "G J C (H S ) Sirivara"
    40 = "SSRHS Pura Turuvekere"
    41 = "G J C NedasaleKunigal"
    43 = "Gubbi"
    44 = "S S R H S Gangalgonte"
    45 = "G H S ArasikerePavagada"
    46 = "G H S YelerampuraKoratagere"
    47 = "G J C Sira"
    48 = "G J C Turuvekere"
    50 = "G H S Ankasandra Gubbi"
    51 = "M G M H S ArasikerePavagadaMadhugiri"
    52 = "G H S HonnashettihalliGubbi"
    53 = "G H S Shivanagere"
    54 = "G H S Thondoti"
    55 = "Sri Swamy Vivekananda High School Kuppur"
    56 = "Sri Siddhartha Rural High School"
    57 = "Dr. Ambedkar Rural High SchoolAralaguppeTiptur"
}

foreach ($row in $addresses.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $addresses[$row]
}
